$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: torta
$ws.Range("C2").Value = "5.0-leche,1.0-vainilla,1.0-huevos,2.0-harinita,"

# C3: kuchen manzana
$ws.Range("C3").Value = "2.0-manzana,5.0-huevos,1.0-harinita,"

# C4: queque
$ws.Range("C4").Value = "1.0-vainilla,2.0-huevos,5.0-harinita,"

# C6: pie de limon
$ws.Range("C6").Value = "5.0-merengue,2.0-limon,1.0-crema,5.0-huevos,4.0-harinita,"
